# Update crypto price/volume data per Jan 26 2024 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''41.327.08'
$ws.Range("E2").Value = '  +3.34%  '
$ws.Range("D3").Value = '''2.252.70'
$ws.Range("E3").Value = '  +1.96%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '''302.82'
$ws.Range("D6").Value = '''91.70'
$ws.Range("E6").Value = '  +4.17%  '
$ws.Range("E7").Value = '  +2.20%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '''0.485'
$ws.Range("E9").Value = '  +2.42%  '
$ws.Range("D10").Value = '''53.93'
$ws.Range("E10").Value = '  +8.03%  '
$ws.Range("D11").Value = '''32.16'
$ws.Range("E11").Value = '  +6.83%  '
$ws.Range("D12").Value = '''0.0794'
$ws.Range("E12").Value = '  +2.09%  '
$ws.Range("E13").Value = '  +2.97%  '
$ws.Range("D14").Value = '''6.61'
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("D15").Value = '''2.599.58'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("E16").Value = '  +2.76%  '
$ws.Range("D17").Value = '''2.255.54'
$ws.Range("E17").Value = '  -5.45%  '
$ws.Range("D18").Value = '''0.750'
$ws.Range("E18").Value = '  +3.37%  '
$ws.Range("D19").Value = '''41.236.04'
$ws.Range("E19").Value = '  +3.31%  '
$ws.Range("D20").Value = '''12.29'
$ws.Range("E20").Value = '  +9.11%  '
$ws.Range("D21").Value = '''0.0₃0904'
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("D22").Value = '''5.90'
$ws.Range("E22").Value = '  +2.25%  '
$ws.Range("D23").Value = '''66.78'
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("D24").Value = '''240.84'
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("E25").Value = '  +4.56%  '
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("D28").Value = '''23.78'
$ws.Range("E28").Value = '  +5.77%  '
$ws.Range("D29").Value = '''2.19'
$ws.Range("E29").Value = '  -1.66%  '
$ws.Range("E30").Value = '  +5.07%  '
$ws.Range("D31").Value = '''158.26'
$ws.Range("E31").Value = '  +1.16%  '
$ws.Range("D32").Value = '''33.54'
$ws.Range("E32").Value = '  +6.99%  '
$ws.Range("D33").Value = '''0.998'
$ws.Range("E33").Value = '  -0.08%  '
$ws.Range("D34").Value = '''5.20'
$ws.Range("E34").Value = '  +6.44%  '
$ws.Range("D35").Value = '''0.0736'
$ws.Range("E36").Value = '  +7.01%  '
$ws.Range("E37").Value = '  +0.99%  '
$ws.Range("D38").Value = '''16.69'
$ws.Range("E38").Value = '  +8.99%  '
$ws.Range("E39").Value = '  +2.99%  '
$ws.Range("D40").Value = '''0.104'
$ws.Range("E40").Value = '  +6.15%  '
$ws.Range("E41").Value = '  +6.30%  '
$ws.Range("E42").Value = '  +5.70%  '
$ws.Range("D43").Value = '''20.64'
$ws.Range("E43").Value = '  +18.59%  '
$ws.Range("D44").Value = '''2.062.70'
$ws.Range("E44").Value = '  -2.77%  '
$ws.Range("E45").Value = '  +3.62%  '
$ws.Range("D46").Value = '''10.23'
$ws.Range("E46").Value = '  +6.46%  '
$ws.Range("D47").Value = '''2.97'
$ws.Range("E47").Value = '  +11.90%  '
$ws.Range("E48").Value = '  -2.94%  '
$ws.Range("D49").Value = '''2.471.31'
$ws.Range("E49").Value = '  +2.06%  '
$ws.Range("D50").Value = '''1.52'
$ws.Range("E50").Value = '  +2.89%  '
$ws.Range("E51").Value = '  +3.77%  '
